$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Plain value edits (style/number format unchanged)
# ---------------------------------------------------------------------------
$ws.Range("H3").Value  = 0.001
$ws.Range("I3").Value  = 0.001
$ws.Range("H4").Value  = 0.001
$ws.Range("I5").Value  = 0.01
$ws.Range("I12").Value = 0.001
$ws.Range("H18").Value = 0.001
$ws.Range("H25").Value = 0.001
$ws.Range("H32").Value = 0.001
$ws.Range("H33").Value = 0.01
$ws.Range("T34").Value = 0

# ---------------------------------------------------------------------------
# H2 / I2: were unstyled (General), now get the "0.000" number format and a
# new value of 0
# ---------------------------------------------------------------------------
$cells2 = @("H2", "I2")
foreach ($addr in $cells2) {
    $c = $ws.Range($addr)
    $c.Value = 0
    $c.NumberFormat = "0.000"
}

# ---------------------------------------------------------------------------
# H/I pairs on the "template" rows: value 0.001 -> 0 and number format
# switches from General (inherited via the row's protected style) to
# "0.000" while keeping the existing protection (locked, not hidden)
# ---------------------------------------------------------------------------
$templateRows = @(6, 11, 13, 19, 24, 26, 31, 34, 40, 45, 48)
foreach ($r in $templateRows) {
    foreach ($col in @("H", "I")) {
        $c = $ws.Range("$col$r")
        $c.Value = 0
        $c.NumberFormat = "0.000"
    }
}

# ---------------------------------------------------------------------------
# A/B/D label cells on those same template rows: pick up the protected
# ("locked / not hidden") style that the rest of the row already carries
# ---------------------------------------------------------------------------
$labelRows = @(11, 13, 19, 24, 26, 31, 34, 40, 45, 48)
foreach ($r in $labelRows) {
    foreach ($col in @("A", "B", "D")) {
        $ws.Range("$col$r").Locked = $true
    }
}
